$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.568.31'
$ws.Range('E2').Value = '  -0.12%  '

# Row 3
$ws.Range('D3').Value = '1.752.19'
$ws.Range('E3').Value = '  -0.28%  '

# Row 4
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
$ws.Range('D5').Value = '''324.35'
$ws.Range('E5').Value = '  -0.13%  '

# Row 6
$ws.Range('E6').Value = '  +0.07%  '

# Row 7
$ws.Range('D7').Value = '''0.4483'
$ws.Range('E7').Value = '  +0.40%  '

# Row 8
$ws.Range('D8').Value = '''0.3563'
$ws.Range('E8').Value = '  -1.53%  '

# Row 9
$ws.Range('D9').Value = '''0.07470'
$ws.Range('E9').Value = '  -1.04%  '

# Row 10
$ws.Range('D10').Value = '''41.46'
$ws.Range('E10').Value = '  -1.49%  '

# Row 11
$ws.Range('D11').Value = '''1.082'
$ws.Range('E11').Value = '  -2.26%  '

# Row 12
$ws.Range('E12').Value = '  +0.07%  '

# Row 13
$ws.Range('D13').Value = '''20.73'
$ws.Range('E13').Value = '  -0.23%  '

# Row 14
$ws.Range('D14').Value = '''5.987'
$ws.Range('E14').Value = '  -1.25%  '

# Row 15
$ws.Range('D15').Value = '''7.151'
$ws.Range('E15').Value = '  -0.97%  '

# Row 16
$ws.Range('D16').Value = '1.755.91'
$ws.Range('E16').Value = '  -0.01%  '

# Row 17
$ws.Range('D17').Value = '''93.71'
$ws.Range('E17').Value = '  +1.06%  '

# Row 18
$ws.Range('D18').Value = '''0.00001057'
$ws.Range('E18').Value = '  -0.86%  '

# Row 19
$ws.Range('D19').Value = '''0.06370'
$ws.Range('E19').Value = '  -0.75%  '

# Row 20
$ws.Range('E20').Value = '  +0.10%  '

# Row 21
$ws.Range('D21').Value = '''17.16'
$ws.Range('E21').Value = '  +0.44%  '

# Row 22
$ws.Range('D22').Value = '''5.738'
$ws.Range('E22').Value = '  -1.88%  '

# Row 23
$ws.Range('D23').Value = '27.620.05'
$ws.Range('E23').Value = '  -0.09%  '

# Row 24
$ws.Range('E24').Value = '  -0.46%  '

# Row 25
$ws.Range('E25').Value = '  -0.61%  '

# Row 26
$ws.Range('D26').Value = '''165.58'
$ws.Range('E26').Value = '  +1.91%  '

# Row 27
$ws.Range('D27').Value = '''20.16'
$ws.Range('E27').Value = '  -1.27%  '

# Row 28
$ws.Range('D28').Value = '1.956.92'
$ws.Range('E28').Value = '  +0.05%  '

# Row 29
$ws.Range('D29').Value = '''2.097'
$ws.Range('E29').Value = '  -1.87%  '

# Row 30
$ws.Range('D30').Value = '''125.64'
$ws.Range('E30').Value = '  -0.18%  '

# Row 31
$ws.Range('D31').Value = '''1.094'
$ws.Range('E31').Value = '  -0.25%  '

# Row 32
$ws.Range('E32').Value = '  +1.77%  '

# Row 33
$ws.Range('D33').Value = '''3.662'
$ws.Range('E33').Value = '  +0.06%  '

# Row 34
$ws.Range('D34').Value = '''5.520'
$ws.Range('E34').Value = '  -0.79%  '

# Row 35
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').Value = '''11.79'
$ws.Range('E35').Value = '  -3.11%  '

# Row 36
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').Value = '''0.02286'
$ws.Range('E36').Value = '  -0.70%  '

# Row 37
$ws.Range('D37').Value = '''0.2095'
$ws.Range('E37').Value = '  -0.84%  '

# Row 38
$ws.Range('D38').Value = '''0.06022'
$ws.Range('E38').Value = '  +0.24%  '

# Row 39
$ws.Range('D39').Value = '''0.6295'
$ws.Range('E39').Value = '  -1.87%  '

# Row 40
$ws.Range('D40').Value = '''4.930'
$ws.Range('E40').Value = '  -0.28%  '

# Row 41
$ws.Range('D41').Value = '''1.181'
$ws.Range('E41').Value = '  -0.39%  '

# Row 42
$ws.Range('D42').Value = '''1.391'
$ws.Range('E42').Value = '  -0.62%  '

# Row 43
$ws.Range('D43').Value = '''7.785'
$ws.Range('E43').Value = '  -1.21%  '

# Row 44
$ws.Range('D44').Value = '''13.19'
$ws.Range('E44').Value = '  -1.27%  '

# Row 45
$ws.Range('E45').Value = '  +0.04%  '

# Row 46
$ws.Range('D46').Value = '''0.5877'
$ws.Range('E46').Value = '  -0.83%  '

# Row 47
$ws.Range('D47').Value = '''122.23'
$ws.Range('E47').Value = '  +0.23%  '

# Row 48
$ws.Range('D48').Value = '''1.936'
$ws.Range('E48').Value = '  -2.29%  '

# Row 49
$ws.Range('D49').Value = '''0.06877'
$ws.Range('E49').Value = '  +0.03%  '

# Row 50
$ws.Range('D50').Value = '''1.130'
$ws.Range('E50').Value = '  -3.21%  '

# Row 51
$ws.Range('D51').Value = '''71.67'
$ws.Range('E51').Value = '  -1.65%  '
